$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "added prior to push to github 5/13/2020 14:30"
$ws.Range("A11").Select() | Out-Null
